$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.144.90"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "2.479.66"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("E4").Value = "  -0.01%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "584.52"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.51%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "172.80"
$c.ClearFormats()
$ws.Range("E6").Value = "  +3.27%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +0.37%  "

$ws.Range("D9").Value = "2.479.28"

$ws.Range("E10").Value = "  +2.82%  "

$ws.Range("E11").Value = "  +1.04%  "

$ws.Range("E12").Value = "  -0.35%  "

$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").Value = "2.932.43"
$ws.Range("E14").Value = "  +0.50%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "25.51"
$c.ClearFormats()
$ws.Range("E15").Value = "  +0.09%  "

$ws.Range("D16").Value = "67.027.81"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("D18").Value = "2.406.74"
$ws.Range("E18").Value = "  -2.00%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.62"
$c.ClearFormats()
$ws.Range("E19").Value = "  +0.23%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "10.95"
$c.ClearFormats()
$ws.Range("E20").Value = "  -3.02%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "350.03"
$c.ClearFormats()
$ws.Range("E21").Value = "  -1.55%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.02"
$c.ClearFormats()
$ws.Range("E22").Value = "  -0.62%  "

$ws.Range("E23").Value = "  +0.08%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "68.94"
$c.ClearFormats()
$ws.Range("E24").Value = "  -0.66%  "

$ws.Range("E25").Value = "  +0.45%  "

$ws.Range("E26").Value = "  +1.90%  "

$ws.Range("E27").Value = "  +0.59%  "

$ws.Range("D28").Value = "2.606.01"
$ws.Range("E28").Value = "  +0.29%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("E30").Value = "  +0.87%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "504.98"
$c.ClearFormats()
$ws.Range("E31").Value = "  -1.06%  "

$ws.Range("E32").Value = "  -0.96%  "

$ws.Range("E33").Value = "  +0.75%  "

$ws.Range("E34").Value = "  -0.87%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E35").Value = "  -0.03%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "161.84"
$c.ClearFormats()
$ws.Range("E36").Value = "  +2.36%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.118"
$c.ClearFormats()
$ws.Range("E37").Value = "  -0.22%  "

$ws.Range("E38").Value = "  +0.57%  "

$ws.Range("E39").Value = "  -1.52%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.33"
$c.ClearFormats()
$ws.Range("E40").Value = "  -0.67%  "

$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("E42").Value = "  +0.75%  "

$ws.Range("E43").Value = "  +1.27%  "

$ws.Range("E44").Value = "  +0.65%  "

$ws.Range("E45").Value = "  +2.37%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "143.21"
$c.ClearFormats()
$ws.Range("E46").Value = "  +1.43%  "

$ws.Range("D47").Value = "0.0₆0266"
$ws.Range("E47").Value = "  +6.19%  "

$ws.Range("E48").Value = "  +0.41%  "

$ws.Range("E49").Value = "  +0.14%  "

$ws.Range("E50").Value = "  +0.67%  "

$ws.Range("E51").Value = "  -0.77%  "
